{"js": "const newValues = [[\"22+67=89\", \"10+75=85\", \"52-37=15\", \"72-65=7\", \"8+2=10\"], [\"41-29=12\", \"12+3=15\", \"13+27=40\", \"86-58=28\", \"62-16=46\"], [\"37-6=31\", \"80-57=23\", \"88-70=18\", \"72-17=55\", \"5+28=33\"], [\"28+56=84\", \"80-73=7\", \"73-26=47\", \"69+3=72\", \"85-12=73\"], [\"85-14=71\", \"28+28=56\", \"99-77=22\", \"15+17=32\", \"71+2=73\"], [\"75-29=46\", \"9+31=40\", \"28+36=64\", \"92-68=24\", \"85-6=79\"], [\"25+70=95\", \"0+86=86\", \"20+79=99\", \"22-13=9\", \"0+71=71\"], [\"79-19=60\", \"5+29=34\", \"24-11=13\", \"45-11=34\", \"20+17=37\"], [\"6+53=59\", \"91-82=9\", \"73-7=66\", \"9+71=80\", \"24+75=99\"], [\"92-60=32\", \"75-65=10\", \"50+25=75\", \"98-80=18\", \"9+64=73\"], [\"4+63=67\", \"27+38=65\", \"24+1=25\", \"54-16=38\", \"86-27=59\"], [\"28-5=23\", \"53+21=74\", \"60-11=49\", \"11-5=6\", \"58-15=43\"], [\"31-27=4\", \"17+42=59\", \"72+25=97\", \"26-19=7\", \"12+37=49\"], [\"78-60=18\", \"2+76=78\", \"62-21=41\", \"73-55=18\", \"27-3=24\"], [\"42-36=6\", \"59-4=55\", \"94-37=57\", \"97-50=47\", \"81-25=56\"], [\"83-65=18\", \"18+32=50\", \"44-6=38\", \"21+13=34\", \"29+11=40\"], [\"11+71=82\", \"81+1=82\", \"90-67=23\", \"39-17=22\", \"71+12=83\"], [\"86-21=65\", \"66-10=56\", \"40+49=89\", \"38+5=43\", \"25+1=26\"], [\"0+82=82\", \"45-42=3\", \"50-42=8\", \"49+31=80\", \"87-51=36\"], [\"72-37=35\", \"78-19=59\", \"18+40=58\", \"79-7=72\", \"62-39=23\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    table.getCell(r, c).value = newValues[r][c];\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"22+67=89\",\"10+75=85\",\"52-37=15\",\"72-65=7\",\"8+2=10\"),\n    @(\"41-29=12\",\"12+3=15\",\"13+27=40\",\"86-58=28\",\"62-16=46\"),\n    @(\"37-6=31\",\"80-57=23\",\"88-70=18\",\"72-17=55\",\"5+28=33\"),\n    @(\"28+56=84\",\"80-73=7\",\"73-26=47\",\"69+3=72\",\"85-12=73\"),\n    @(\"85-14=71\",\"28+28=56\",\"99-77=22\",\"15+17=32\",\"71+2=73\"),\n    @(\"75-29=46\",\"9+31=40\",\"28+36=64\",\"92-68=24\",\"85-6=79\"),\n    @(\"25+70=95\",\"0+86=86\",\"20+79=99\",\"22-13=9\",\"0+71=71\"),\n    @(\"79-19=60\",\"5+29=34\",\"24-11=13\",\"45-11=34\",\"20+17=37\"),\n    @(\"6+53=59\",\"91-82=9\",\"73-7=66\",\"9+71=80\",\"24+75=99\"),\n    @(\"92-60=32\",\"75-65=10\",\"50+25=75\",\"98-80=18\",\"9+64=73\"),\n    @(\"4+63=67\",\"27+38=65\",\"24+1=25\",\"54-16=38\",\"86-27=59\"),\n    @(\"28-5=23\",\"53+21=74\",\"60-11=49\",\"11-5=6\",\"58-15=43\"),\n    @(\"31-27=4\",\"17+42=59\",\"72+25=97\",\"26-19=7\",\"12+37=49\"),\n    @(\"78-60=18\",\"2+76=78\",\"62-21=41\",\"73-55=18\",\"27-3=24\"),\n    @(\"42-36=6\",\"59-4=55\",\"94-37=57\",\"97-50=47\",\"81-25=56\"),\n    @(\"83-65=18\",\"18+32=50\",\"44-6=38\",\"21+13=34\",\"29+11=40\"),\n    @(\"11+71=82\",\"81+1=82\",\"90-67=23\",\"39-17=22\",\"71+12=83\"),\n    @(\"86-21=65\",\"66-10=56\",\"40+49=89\",\"38+5=43\",\"25+1=26\"),\n    @(\"0+82=82\",\"45-42=3\",\"50-42=8\",\"49+31=80\",\"87-51=36\"),\n    @(\"72-37=35\",\"78-19=59\",\"18+40=58\",\"79-7=72\",\"62-39=23\")\n)\n\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n    $row = $newValues[$r]\n    for ($c = 0; $c -lt $row.Count; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $row[$c]\n    }\n}\n"}
